$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" (F1) and "is_enabled" (G1) header columns.
# This shifts the remaining "rem" column (was H1) left to F1,
# matching the diff which deletes the two corresponding shared-string
# entries and the two corresponding worksheet cells.
$ws.Range("F1:G1").EntireColumn.Delete()
